$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Exchange") - shifts Symbol..Total P/L right by one.
$ws.Columns("B:B").Insert()

# Clear the formatting the Insert copied into the new data cell (B2 should be unstyled,
# matching every other non-header data cell).
$ws.Cells.Item(2, 2).ClearFormats()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 2).Value2 = "Exchange"
$ws.Cells.Item(1, 12).Value2 = "Strategy"

# --- Data row (row 2) ---
$ws.Cells.Item(2, 2).Value2 = "ByBit"

# "From" date needs to stay a plain text string (not get auto-parsed into a date serial).
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value2 = "2021-10-01"
$ws.Cells.Item(2, 4).ClearFormats()

$ws.Cells.Item(2, 8).Value2 = 6
$ws.Cells.Item(2, 9).Value2 = 4
$ws.Cells.Item(2, 12).Value2 = "MACD"
$ws.Cells.Item(2, 13).Value2 = 7
$ws.Cells.Item(2, 14).Value2 = 6
$ws.Cells.Item(2, 15).Value2 = 13

# "Success Rate" needs to stay a plain text string (not get auto-parsed into a percentage).
$ws.Cells.Item(2, 16).NumberFormat = "@"
$ws.Cells.Item(2, 16).Value2 = "53.8%"
$ws.Cells.Item(2, 16).ClearFormats()

$ws.Cells.Item(2, 17).Value2 = 0
$ws.Cells.Item(2, 18).Value2 = 4
$ws.Cells.Item(2, 19).Value2 = 4200
$ws.Cells.Item(2, 20).Value2 = -2400
$ws.Cells.Item(2, 21).Value2 = 132.05
$ws.Cells.Item(2, 22).Value2 = 1667.95
